$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.026490916865734
$ws.Range("E2").Value = 1.035281354641554
$ws.Range("F2").Value = 1.043026249814855
$ws.Range("J2").Value = 1.031654064059758
$ws.Range("L2").Value = 1.03807840888099
$ws.Range("M2").Value = 1.045801277101354
$ws.Range("N2").Value = 1.033119131163119
# Row 3
$ws.Range("C3").Value = 1.027842156397847
$ws.Range("E3").Value = 1.036523650962609
$ws.Range("F3").Value = 1.044448464980237
$ws.Range("J3").Value = 1.032643018763664
$ws.Range("L3").Value = 1.039128916684576
$ws.Range("M3").Value = 1.047032844707646
$ws.Range("N3").Value = 1.034109490296118
# Row 4
$ws.Range("C4").Value = 1.028716173043294
$ws.Range("E4").Value = 1.037327489195902
$ws.Range("F4").Value = 1.045369012804246
$ws.Range("J4").Value = 1.033282201577143
$ws.Range("L4").Value = 1.039808101438637
$ws.Range("M4").Value = 1.047829512859633
$ws.Range("N4").Value = 1.034749580822508
# Row 5
$ws.Range("C5").Value = 1.029083536151651
$ws.Range("E5").Value = 1.037665423616462
$ws.Range("F5").Value = 1.045756082122717
$ws.Range("J5").Value = 1.033550740576032
$ws.Range("L5").Value = 1.040093498731203
$ws.Range("M5").Value = 1.048164378682872
$ws.Range("N5").Value = 1.035018501177577
# Row 6
$ws.Range("C6").Value = 1.029145213854658
$ws.Range("E6").Value = 1.03772216443854
$ws.Range("F6").Value = 1.045821077047702
$ws.Range("J6").Value = 1.033595819370691
$ws.Range("L6").Value = 1.040141410541831
$ws.Range("M6").Value = 1.048220601097659
$ws.Range("N6").Value = 1.035063643989295
# Row 7
$ws.Range("C7").Value = 1.028721082050322
$ws.Range("E7").Value = 1.037332004684756
$ws.Range("F7").Value = 1.045374184559297
$ws.Range("J7").Value = 1.033285790486255
$ws.Range("L7").Value = 1.039811915446028
$ws.Range("M7").Value = 1.047833987558758
$ws.Range("N7").Value = 1.034753174828282
# Row 8
$ws.Range("C8").Value = 1.026947644726116
$ws.Range("E8").Value = 1.03570119870601
$ws.Range("F8").Value = 1.043506838681758
$ws.Range("J8").Value = 1.031988440333347
$ws.Range("L8").Value = 1.038433551539637
$ws.Range("M8").Value = 1.046217542688354
$ws.Range("N8").Value = 1.033453982289349
# Row 9
$ws.Range("C9").Value = 1.023819926531252
$ws.Range("E9").Value = 1.032827278913259
$ws.Range("F9").Value = 1.04021829376622
$ws.Range("J9").Value = 1.029696565096051
$ws.Range("L9").Value = 1.036000241370181
$ws.Range("M9").Value = 1.04336717205273
$ws.Range("N9").Value = 1.031158852326335
# Row 10
$ws.Range("C10").Value = 1.021732687406676
$ws.Range("E10").Value = 1.030910964937144
$ws.Range("F10").Value = 1.038026988218159
$ws.Range("J10").Value = 1.028164574574658
$ws.Range("L10").Value = 1.034374845117875
$ws.Range("M10").Value = 1.041465360200861
$ws.Range("N10").Value = 1.029624686202679
# Row 11
$ws.Range("C11").Value = 1.020828330211808
$ws.Range("E11").Value = 1.030081045017767
$ws.Range("F11").Value = 1.037078321730686
$ws.Range("J11").Value = 1.027500197155752
$ws.Range("L11").Value = 1.033670230717229
$ws.Range("M11").Value = 1.040641432365732
$ws.Range("N11").Value = 1.028959365291633
# Row 12
$ws.Range("C12").Value = 1.020492320724776
$ws.Range("E12").Value = 1.029772750100411
$ws.Range("F12").Value = 1.036725967244162
$ws.Range("J12").Value = 1.027253261528089
$ws.Range("L12").Value = 1.033408380493541
$ws.Range("M12").Value = 1.040335319692193
$ws.Range("N12").Value = 1.028712078987055
# Row 13
$ws.Range("C13").Value = 1.020564400104801
$ws.Range("E13").Value = 1.029838881637244
$ws.Range("F13").Value = 1.036801547541352
$ws.Range("J13").Value = 1.027306237190017
$ws.Range("L13").Value = 1.033464553975997
$ws.Range("M13").Value = 1.040400985102833
$ws.Range("N13").Value = 1.028765129880499
# Row 14
$ws.Range("C14").Value = 1.020800557450149
$ws.Range("E14").Value = 1.030055561812355
$ws.Range("F14").Value = 1.037049195569476
$ws.Range("J14").Value = 1.027479788589223
$ws.Range("L14").Value = 1.033648588669356
$ws.Range("M14").Value = 1.040616130434051
$ws.Range("N14").Value = 1.028938927742598
# Row 15
$ws.Range("C15").Value = 1.020946049567315
$ws.Range("E15").Value = 1.030189062087951
$ws.Range("F15").Value = 1.037201782521286
$ws.Range("J15").Value = 1.027586698515009
$ws.Range("L15").Value = 1.033761961831765
$ws.Range("M15").Value = 1.040748679270769
$ws.Range("N15").Value = 1.029045989492742
# Row 16
$ws.Range("C16").Value = 1.021792694161462
$ws.Range("E16").Value = 1.030966040613384
$ws.Range("F16").Value = 1.038089951519511
$ws.Range("J16").Value = 1.028208645405253
$ws.Range("L16").Value = 1.034421590726031
$ws.Range("M16").Value = 1.041520032120773
$ws.Range("N16").Value = 1.029668819618909
# Row 17
$ws.Range("C17").Value = 1.022323615771932
$ws.Range("E17").Value = 1.031453377401598
$ws.Range("F17").Value = 1.038647121976613
$ws.Range("J17").Value = 1.028598501782262
$ws.Range("L17").Value = 1.034835139739721
$ws.Range("M17").Value = 1.042003762906311
$ws.Range("N17").Value = 1.03005922963669
# Row 18
$ws.Range("C18").Value = 1.022633238848841
$ws.Range("E18").Value = 1.031737619160359
$ws.Range("F18").Value = 1.038972128104003
$ws.Range("J18").Value = 1.028825800691129
$ws.Range("L18").Value = 1.035076278268112
$ws.Range("M18").Value = 1.042285873631966
$ws.Range("N18").Value = 1.03028685133608
# Row 19
$ws.Range("C19").Value = 1.02273880323527
$ws.Range("E19").Value = 1.031834536082859
$ws.Range("F19").Value = 1.039082950004158
$ws.Range("J19").Value = 1.028903287282061
$ws.Range("L19").Value = 1.035158487224475
$ws.Range("M19").Value = 1.042382059235768
$ws.Range("N19").Value = 1.030364447966858
# Row 20
$ws.Range("C20").Value = 1.022266658587826
$ws.Range("E20").Value = 1.031401092208853
$ws.Range("F20").Value = 1.03858734103306
$ws.Range("J20").Value = 1.028556684018271
$ws.Range("L20").Value = 1.034790777901109
$ws.Range("M20").Value = 1.041951867472193
$ws.Range("N20").Value = 1.030017352486678
# Row 21
$ws.Range("C21").Value = 1.020731017553108
$ws.Range("E21").Value = 1.029991755644657
$ws.Range("F21").Value = 1.036976268817576
$ws.Range("J21").Value = 1.02742868637165
$ws.Range("L21").Value = 1.033594398522078
$ws.Range("M21").Value = 1.040552777442377
$ws.Range("N21").Value = 1.028887752954015
# Row 22
$ws.Range("C22").Value = 1.019764970093796
$ws.Range("E22").Value = 1.029105498996335
$ws.Range("F22").Value = 1.035963451356113
$ws.Range("J22").Value = 1.026718563481539
$ws.Range("L22").Value = 1.032841461854166
$ws.Range("M22").Value = 1.039672711468771
$ws.Range("N22").Value = 1.028176621607948
# Row 23
$ws.Range("C23").Value = 1.020277142110416
$ws.Range("E23").Value = 1.029575336107435
$ws.Range("F23").Value = 1.036500354652359
$ws.Range("J23").Value = 1.027095100184893
$ws.Range("L23").Value = 1.033240677867889
$ws.Range("M23").Value = 1.040139290747095
$ws.Range("N23").Value = 1.02855369303661
# Row 24
$ws.Range("C24").Value = 1.022292395255737
$ws.Range("E24").Value = 1.031424717678506
$ws.Range("F24").Value = 1.038614353412235
$ws.Range("J24").Value = 1.028575579967525
$ws.Range("L24").Value = 1.034810823346056
$ws.Range("M24").Value = 1.041975316912289
$ws.Range("N24").Value = 1.030036275270347
# Row 25
$ws.Range("C25").Value = 1.024628864977666
$ws.Range("E25").Value = 1.033570305248744
$ws.Range("F25").Value = 1.041068256048337
$ws.Range("J25").Value = 1.030289773934296
$ws.Range("L25").Value = 1.036629858356136
$ws.Range("M25").Value = 1.044104320226558
$ws.Range("N25").Value = 1.031752903589172
